$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('realm; field', '分野|ぶんや'),
    @('Mr./Ms. Ono', '小野さん|おのさん'),
    @('Nagano', '長野|ながの'),
    @('baseball', '野球|やきゅう'),
    @('vegetable', '野菜|やさい'),
    @('to learn', '習う|ならう'),
    @('calligraphy', '習字|しゅうじ'),
    @('practice', '練習|れんしゅう'),
    @('habit; custom', '習慣|しゅうかん'),
    @('mainly', '主に|おもに'),
    @('husband', 'ご主人|ごしゅじん'),
    @('housewife', '主婦|しゅふ'),
    @('subject of a sentence', '主語|しゅご'),
    @('owner', '持ち主|もちぬし'),
    @('twenty-five years old', '二十五歳|にじゅうごさい'),
    @('year-end gift', 'お歳暮|おせいぼ'),
    @('twenty years old', '二十歳|はたち'),
    @('to collect', '集める|あつめる'),
    @('(magazine) feature', '特集|とくしゅう'),
    @('to concentrate', '集中する|しゅうちゅうする'),
    @('to make public; to give a presentation', '発表する|はっぴょうする'),
    @('pronunciation', '発音|はつおん'),
    @('departure', '出発|しゅっぱつ'),
    @('invention', '発明|はつめい'),
    @('to make public; to give a presentation', '発表する|はっぴょうする'),
    @('cover page', '表紙|ひょうし'),
    @('expression', '表現|ひょうげん'),
    @('to express', '表す|あらわす'),
    @('the front', '表|おもて'),
    @('work (of art, etc.)', '作品|さくひん'),
    @('elegant', '上品な|じょうひんな'),
    @('magic', '手品|てじな'),
    @('article of merchandise', '品物|しなもの'),
    @('photograph', '写真|しゃしん'),
    @('to copy', '写す|うつす'),
    @('sketch', '写生|しゃせい'),
    @('to describe', '描写する|びょうしゃする'),
    @('photograph', '写真|しゃしん'),
    @('center', '真ん中|まんなか'),
    @('right above...', '～の真上|～のまうえ'),
    @('midnight', '真夜中|まよなか'),
    @('character', '文字|もじ'),
    @('deficit', '赤字|あかじ'),
    @('family name', '名字|みょうじ'),
    @('uppercase letters', '大文字|おおもじ'),
    @('Chinese character', '漢字|かんじ'),
    @('activity', '活動|かつどう'),
    @('life; living', '生活|せいかつ'),
    @('active', '活発な|かっぱつな'),
    @('conjugation', '活用|かつよう'),
    @('liveliness', '活気|かっき'),
    @('to get married', '結婚する|けっこんする'),
    @('result', '結果|けっか'),
    @('conclusion', '結論|けつろん'),
    @('to tie a knot', '結ぶ|むすぶ'),
    @('divorce', '離婚|りこん'),
    @('fiancé(e)', '婚約者|こんやくしゃ'),
    @('newlywed', '新婚|しんこん'),
    @('unmarried', '未婚|みこん'),
    @('to walk', '歩く|あるく'),
    @('sidewalk', '歩道|ほどう'),
    @('to stroll', '散歩する|さんぽする'),
    @('one step', '一歩|いっぽ'),
    @('progress', '進歩|しんぽ')
)

$startRow = 85
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
